# Fix the "pyhton" -> "python" typo on slide 7 ("22.2 Upgrade Database"),
# in the subtitle placeholder's last bullet:
#   "> pyhton manage.py migrate "  ->  "> python manage.py migrate "
#
# The original run layout for that paragraph was:
#   Run1 "> "      Run2 "pyhton" (err="1")      Run3 " manage.py migrate "
# After the fix it becomes two runs:
#   Run1 "> python "                           Run2 "manage.py migrate "

$p = $ppt.ActivePresentation

$targetSlideIndex = 7
$targetShapeIndex = 3
$needle = "pyhton"
$oldFragment = "> pyhton "
$newFragment = "> python "

$slide = $p.Slides.Item($targetSlideIndex)
$shape = $slide.Shapes.Item($targetShapeIndex)
$textRange = $shape.TextFrame.TextRange

$paragraphCount = $textRange.Paragraphs().Count
$fixed = $false

for ($i = 1; $i -le $paragraphCount; $i++) {
    $para = $textRange.Paragraphs($i, 1)
    if ($para.Text.IndexOf($needle) -ge 0) {
        $paraStart = $para.Start
        $offset = $para.Text.IndexOf($oldFragment)
        if ($offset -ge 0) {
            $fragmentRange = $textRange.Characters($paraStart + $offset, $oldFragment.Length)
            $fragmentRange.Text = $newFragment
            $fixed = $true
        }
    }
}

if (-not $fixed) {
    throw "Could not locate the '$needle' typo to fix."
}
